# [맘파] StageTable StageTimer 칼럼 추가
# Adds a new "StageTimer" column (data-dictionary row on "#Index" sheet,
# and the actual data column L on the "stage" sheet).

$wb = $excel.ActiveWorkbook

$wsIndex = $wb.Worksheets.Item("#Index")
$wsStage = $wb.Worksheets.Item("stage")

# ---------------------------------------------------------------------
# "#Index" sheet: document the new column in the data-dictionary table
# (row 18, right under the existing "StageFile" row).
# ---------------------------------------------------------------------
$wsIndex.Range("A18").Value = "StageTimer"
$wsIndex.Range("B18").Value = "int"
$wsIndex.Range("C18").Value = "스테이지 제한 시간 [ 초 단위로 설정 ]"

# ---------------------------------------------------------------------
# "stage" sheet: add the StageTimer column (L) with header/type rows and
# per-stage values.
# ---------------------------------------------------------------------
$wsStage.Range("L1").Value = "StageTimer"
$wsStage.Range("L2").Value = "int"

$wsStage.Range("L3").Value = 120
$wsStage.Range("L4").Value = 150
$wsStage.Range("L5").Value = 240
$wsStage.Range("L6").Value = 240
$wsStage.Range("L7").Value = 420
$wsStage.Range("L8").Value = 210
$wsStage.Range("L9").Value = 210
$wsStage.Range("L10").Value = 210
$wsStage.Range("L11").Value = 300
$wsStage.Range("L12").Value = 420

# Column width tweaks that came along with the new column:
#  - column B (Name) widened a bit
#  - column L (new StageTimer) gets its own custom width
$wsStage.Columns.Item(2).ColumnWidth = 20.334
$wsStage.Columns.Item(12).ColumnWidth = 9.834
